$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Điều 1 heading paragraph (bookmark "dieu_1"): the original text was
#    split across many runs, each individually wrapped in <w:proofErr
#    spellStart/spellEnd> pairs left over from spell-checking. Re-inserting
#    the paragraph's own WordOpenXML collapses the runs into a single run
#    and drops the stray proofErr markers, while keeping the bookmark and
#    bold formatting intact.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML($p1.Range.WordOpenXML)

# ---------------------------------------------------------------------------
# 2) "1. Sửa đổi, bổ sung một số khoản của Điều 3 như sau:" paragraph: same
#    proofErr/run-splitting cleanup, while keeping the "dc_1" bookmark around
#    "Điều 3".
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML($p2.Range.WordOpenXML)

# ---------------------------------------------------------------------------
# 3) "a) Sửa đổi, bổ sung khoản 4 như sau:" paragraph: same cleanup, keeping
#    the "diem_a_1_1" bookmark around "a) Sửa đổi, bổ sung".
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML($p3.Range.WordOpenXML)

# ---------------------------------------------------------------------------
# 4) Drop the closing "approved by the National Assembly" paragraphs and the
#    signature table, then add a trailing space to the "hiệu lực thi hành"
#    sentence so it reads "...năm 2025. " immediately before the sectPr.
# ---------------------------------------------------------------------------
$d.Tables(1).Delete()

$trailingPara = $d.Paragraphs($d.Paragraphs.Count)
while ($trailingPara.Range.Text.Trim().Length -eq 0 -and $d.Tables.Count -eq 0) {
    $prevText = $d.Paragraphs($d.Paragraphs.Count - 1).Range.Text
    if ($prevText.StartsWith("Luật này có hiệu lực thi hành")) {
        break
    }
    $trailingPara.Range.Delete()
    $trailingPara = $d.Paragraphs($d.Paragraphs.Count)
}

# Remove the italic "Luật này được Quốc hội ... thông qua ngày ..." paragraph
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
if ($lastPara.Range.Text.StartsWith("Luật này được Quốc hội")) {
    $lastPara.Range.Delete()
}

$effectivePara = $d.Paragraphs($d.Paragraphs.Count)
$r = $effectivePara.Range
$insPoint = $d.Range($r.End - 1, $r.End - 1)
$insPoint.InsertAfter(" ")
